$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WR")

# Week 16 log: B.Perriman is no longer on the roster - remove his row.
$ws.Rows.Item(5).Delete()

# Season sim from Week 17 added a new player to the WR roster.
$ws.Range("A8").Value = "D.Newsome"
for ($c = 2; $c -le 10; $c++) {
    $ws.Cells.Item(8, $c).Value = 0
}

# Reflect the WR sheet as the active sheet/selection, matching the saved view.
$ws.Activate()
$ws.Range("A5:J8").Select()
